# Apply the "TRANSFER" sheet insertion + TEMPORARY_TRANSFER.TRANSFER_ID column edit
$wb = $excel.ActiveWorkbook

# --- 1. Insert new "TRANSFER" sheet between ACCOUNT and REGULAR_TRANSFER ---
$regularTransferSheet = $wb.Worksheets.Item("REGULAR_TRANSFER")
$newSheet = $wb.Worksheets.Add($regularTransferSheet)
$newSheet.Name = "TRANSFER"

# Header row
$newSheet.Range("A1").Value = "ID"
$newSheet.Range("B1").Value = "TITLE"
$newSheet.Range("C1").Value = "USER_ID"

# Data rows
$newSheet.Range("A2").Value = 1
$newSheet.Range("B2").Value = "2023/06"
$newSheet.Range("C2").Value = 1

$newSheet.Range("A3").Value = 2
$newSheet.Range("B3").Value = "2023/07"
$newSheet.Range("C3").Value = 1

$newSheet.Range("A4").Value = 3
$newSheet.Range("B4").Value = "2022/09"
$newSheet.Range("C4").Value = 2

# --- 2. Add TRANSFER_ID column to TEMPORARY_TRANSFER sheet ---
$tempTransferSheet = $wb.Worksheets.Item("TEMPORARY_TRANSFER")
$tempTransferSheet.Range("G1").Value = "TRANSFER_ID"
$tempTransferSheet.Range("G2").Value = 1
$tempTransferSheet.Range("G3").Value = 2
$tempTransferSheet.Range("G4").Value = 3

# --- 3. Set active sheet to TEMPORARY_TRANSFER ---
$tempTransferSheet.Activate()
